# ---------------------------------------------------------------------------
# MousetuaryTODO.xlsx -- "Update changelog and close off 1.0.1."
#
# The worksheet layout is heavily reshuffled between the two revisions (rows
# inserted, removed and renumbered throughout), so rather than trying to
# patch individual cells in place, the whole used range is cleared first and
# then rebuilt row-by-row to match the new TODO list exactly.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Wipe rows 1-37 completely (values, styles, row heights) so nothing from the
# old layout (e.g. stale row heights) lingers once the grid is rebuilt below.
$ws.Rows("1:37").Delete() | Out-Null

# --- Row 1 ---------------------------------------------------------------
$ws.Range("A1").Value = "Mousetuary TODO list"
$rw = $ws.Range("A1:C1")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = "1.0.1"
$rw = $ws.Range("A2:C2")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = "Feature/bug"
$ws.Range("B3").Value = "Status"
$ws.Range("C3").Value = "Comments/Issues"
$rw = $ws.Range("A3:D3")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 4 ---------------------------------------------------------------
$rw = $ws.Range("A4:C4")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = "Home Screen and widget lists"
$rw = $ws.Range("A5")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 6 ---------------------------------------------------------------
$ws.Rows("6:6").RowHeight = 28.8
$ws.Range("A6").Value = "Sysmenu: minimise, window, exit"
$ws.Range("B6").Value = "Done"
$ws.Range("C6").Value = "Need better handling of maximise from Win titlebar. Bug reported. Once fixed, these buttons can be hidden when windowed."
$rw = $ws.Range("A6:C6")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = "Back and Home to topbar."
$ws.Range("B7").Value = "Done"
$rw = $ws.Range("A7:B7")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 8 ---------------------------------------------------------------
$ws.Range("A8").Value = "Tidy up text and clock in topbar"
$ws.Range("B8").Value = "Done"
$ws.Range("C8").Value = "Smaller and less shouty text in general."
$rw = $ws.Range("A8:C8")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 9 ---------------------------------------------------------------
$ws.Rows("9:9").RowHeight = 28.8
$ws.Range("A9").Value = "Remove silly option thing (sidemenu)"
$ws.Range("B9").Value = "Reinstated"
$ws.Range("C9").Value = "Fixed losing focus by using a window property, like Aeon Nox 5 does. Options button up in top bar. Dismiss by a back button."
$rw = $ws.Range("A9:C9")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 10 ---------------------------------------------------------------
$ws.Rows("10:10").RowHeight = 28.8
$ws.Range("A10").Value = "Sensible defaults"
$ws.Range("B10").Value = "Low priority"
$ws.Range("C10").Value = "Can you control these on a skin-by-skin basis? (like order and enabling of buttons on main menu)"
$rw = $ws.Range("A10:C10")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 11 ---------------------------------------------------------------
$ws.Rows("11:11").RowHeight = 43.2
$ws.Range("A11").Value = "Scroll bars in widget lists"
$ws.Range("B11").Value = "Done"
$ws.Range("C11").Value = "Gives a standard way of scrolling them, since autoscroll is not working consistently. The mouse wheel can be used in the scroll bar (while in the widget rows, it scrolls them horizontally)"
$rw = $ws.Range("A11:D11")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 12 ---------------------------------------------------------------
$rw = $ws.Range("A12")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 13 ---------------------------------------------------------------
$ws.Range("A13").Value = "Video OSD"
$rw = $ws.Range("A13")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 14 ---------------------------------------------------------------
$ws.Rows("14:14").RowHeight = 28.8
$ws.Range("A14").Value = "OSD pause button flickers when mouse moved"
$ws.Range("B14").Value = "Done"
$ws.Range("C14").Value = "Removed DefaultControl setting - seems to have no ill effects"
$rw = $ws.Range("A14:C14")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 15 ---------------------------------------------------------------
$ws.Range("A15").Value = "OSD FF/REW and volume buttons"
$ws.Range("B15").Value = "Done"
$rw = $ws.Range("A15:B15")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 16 ---------------------------------------------------------------
$ws.Range("A16").Value = "OSD skip fwd/back buttons"
$ws.Range("B16").Value = "Done"
$ws.Range("C16").Value = "These do Next/Prev for DVD and SmallStepForward/Back for Live TV"
$rw = $ws.Range("A16:C16")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 17 ---------------------------------------------------------------
$ws.Rows("17:17").RowHeight = 28.8
$ws.Range("A17").Value = "OSD transparent background colour nicer"
$ws.Range("B17").Value = "Done"
$ws.Range("C17").Value = "Make it blue like WMC. Also make a semitransparent blue background for the guide when live video is in progress."
$rw = $ws.Range("A17:C17")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 18 ---------------------------------------------------------------
$ws.Rows("18:18").RowHeight = 28.8
$ws.Range("A18").Value = "Get rid of fwd/back buttons from OSD slider nib"
$ws.Range("B18").Value = "Done"
$rw = $ws.Range("A18:B18")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 20 ---------------------------------------------------------------
$ws.Range("A20").Value = "Get rid of misleading message for slider"
$ws.Range("B20").Value = "Done"
$rw = $ws.Range("A20:B20")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 21 ---------------------------------------------------------------
$ws.Rows("21:21").RowHeight = 28.8
$ws.Range("A21").Value = "Use full guide rather than the abbreviated OSD guide. Remove OSD channel button."
$ws.Range("B21").Value = "Done"
$ws.Range("C21").Value = "Ctrl-g for guide with support of a script (goes with optional keymap)"
$rw = $ws.Range("A21:C21")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 22 ---------------------------------------------------------------
$rw = $ws.Range("A22")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 23 ---------------------------------------------------------------
$ws.Range("A23").Value = "Live TV widget list"
$rw = $ws.Range("A23")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 24 ---------------------------------------------------------------
$ws.Rows("24:24").RowHeight = 28.8
$ws.Range("A24").Value = "Channels and Guide seem redundant."
$ws.Range("B24").Value = "Done-ish"
$ws.Range("C24").Value = "Can we just have the Guide and get rid of Channels? For the moment, just make guide the first widget in the line."
$rw = $ws.Range("A24:C24")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 25 ---------------------------------------------------------------
$ws.Range("A25").Value = "Make Guide come first in widget list"
$ws.Range("B25").Value = "Done"
$ws.Range("C25").Value = "Pending decision on the above"
$rw = $ws.Range("A25:C25")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 27 ---------------------------------------------------------------
$ws.Range("A27").Value = "TV Guide"
$rw = $ws.Range("A27")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 28 ---------------------------------------------------------------
$ws.Rows("28:28").RowHeight = 28.8
$ws.Range("A28").Value = "Guide should roll up and down with the mouse wheel, not left to right. "
$ws.Range("B28").Value = "Done"
$ws.Range("C28").Value = "Done using a keymap. "
$rw = $ws.Range("A28:C28")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 29 ---------------------------------------------------------------
$ws.Rows("29:29").RowHeight = 28.8
$ws.Range("A29").Value = "Guide needs up/down/left/right buttons that appear with a mouse hover, like WMC"
$ws.Range("B29").Value = "Done"
$rw = $ws.Range("A29:B29")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 30 ---------------------------------------------------------------
$ws.Rows("30:30").RowHeight = 28.8
$ws.Range("A30").Value = "Guide rows need to be a little bigger, also guide font"
$ws.Range("B30").Value = "Done"
$rw = $ws.Range("A30:B30")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 32 ---------------------------------------------------------------
$ws.Range("A32").Value = "1.0.2"
$rw = $ws.Range("A32")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 33 ---------------------------------------------------------------
$ws.Range("A33").Value = "Incomplete:"
$rw = $ws.Range("A33")
$rw.Font.Bold = $true
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 34 ---------------------------------------------------------------
$ws.Rows("34:34").RowHeight = 28.8
$ws.Range("A34").Value = "Fix over-aggressive autoscroll in Addons/MyAddons list"
$ws.Range("B34").Value = "Not done"
$ws.Range("C34").Value = "It's awful! It should either autoscroll properly, or have a scroll bar."
$rw = $ws.Range("A34:C34")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 35 ---------------------------------------------------------------
$ws.Rows("35:35").RowHeight = 28.8
$ws.Range("A35").Value = "VideoOSDHelpTextVar for new buttons (variables.xml)"
$ws.Range("B35").Value = "Started"
$ws.Range("C35").Value = "Display left-side ones on left, if possible. "
$rw = $ws.Range("A35:C35")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 36 ---------------------------------------------------------------
$ws.Rows("36:36").RowHeight = 28.8
$ws.Range("A36").Value = "Right-click option to delete recent channels from list"
$ws.Range("B36").Value = "Not skinnable"
$ws.Range("C36").Value = "There is very limited scope for skins to do this (existing context menus ca be added to, but new ones can't be created)"
$rw = $ws.Range("A36:C36")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Row 37 ---------------------------------------------------------------
$ws.Rows("37:37").RowHeight = 28.8
$ws.Range("A37").Value = "OSD slider can report timeshift buffer extent to scale, like WMC does"
$ws.Range("B37").Value = "Not skinnable"
$ws.Range("C37").Value = "Need to display a calculated progress based on several infotags. Can't do this in the skinning engine."
$rw = $ws.Range("A37:C37")
$rw.Font.Bold = $false
$rw.WrapText = $true
$rw.VerticalAlignment = -4160  # xlTop

# --- Sheet view -------------------------------------------------------------
# Target view: no frozen/scrolled topLeftCell; whole row 12 selected
# (activeCell A12, sqref A12:XFD12).
$ws.Rows("12:12").EntireRow.Select() | Out-Null

